$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Count value changes from 6 to 4 (kept as text, matching the
# existing "number stored as text" convention used throughout this sheet)
$ws.Range("B3").Value = "'4"

# Row 8: "tja" / greeting -> timestamp field
$ws.Range("A8").Value = "timestamp"
$ws.Range("B8").Value = "2023-11-23 14:40:00.000"

# Row 9: "hallo" / greeting -> sentfrom field
$ws.Range("A9").Value = "sentfrom"
$ws.Range("B9").Value = "holla"

# Row 10: "bjeff" / bark -> mottaker field
$ws.Range("A10").Value = "mottaker"
$ws.Range("B10").Value = "The Stig"

# New row 11: tema
$ws.Range("A11").Value = "tema"
$ws.Range("B11").Value = "fast car ride"

# New row 12: payload (JSON blob)
$ws.Range("A12").Value = "payload"
$ws.Range("B12").Value = '{"id": 123, "name":"Per Spellmann","address":"alle de som æ ø å Æ Ø Å spelle kan de hører hjemme i spelleland", "wish":"han hadde ei einaste ku"}'
